$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.958.38'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +0.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.950.32'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9988'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.12'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9977'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4877'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +1.35%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2969'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.89%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06828'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.17'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.33'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -4.19%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.942.03'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.93%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07731'
$ws.Range('D13').ClearFormats()
$ws.Range('E14').Value = '  -0.13%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7061'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '282.52'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -4.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '30.974.61'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('E18').Value = '  +0.20%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007724'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.203.04'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.54%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9974'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.500'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -2.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9992'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.20%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.494'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.838'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.20'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.97'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.215'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.1053'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -3.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.410'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.05%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.584'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.70%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.568'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.467'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.92%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04950'
$ws.Range('D34').ClearFormats()
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7660'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.171'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.726'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.10%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02020'
$ws.Range('D38').ClearFormats()
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.695'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.20%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.538'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +9.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.156'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.80'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +7.42%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4494'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.87%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '109.47'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8835'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.24%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.193'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +10.84%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9975'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '981.02'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +7.66%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.405'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.55%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.1263'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.77'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.44%  '
